$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update last_edited_time (column D) for rows 4, 5, 6, 8, 12, 13
$ws.Range("D4").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("D5").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("D6").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("D8").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("D12").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("D13").Value = "2024-07-08T01:58:00.000Z"

# Update numeric report figures for row 13 ("Thang 7" / co so report)
$ws.Range("T13").Value = 3500000
$ws.Range("W13").Value = 8284000
$ws.Range("AA13").Value = 50646000
$ws.Range("AE13").Value = 58930000
$ws.Range("AH13").Value = 41430000
$ws.Range("AK13").Value = 9
$ws.Range("AN13").Value = 17500000
$ws.Range("AQ13").Value = 44930000
